$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: password value changed from numeric 19373469 to the text string "1"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1"
$ws.Range("C2").Style = "Normal"

# A4: value "3" changed from text to a genuine number
$ws.Range("A4").Value = 3
